# Fill in the "Variable type" column of the ER-diagram tables on slide 2
# with hardcoded SQL datatypes, and drop the now-redundant "FK address_ID"
# row from the Orders table.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# --- Table 4 ("Customers Table") -------------------------------------
$customers = $s.Shapes.Item(1).Table
$customers.Cell(3, 2).Shape.TextFrame.TextRange.Text = "Int"          # PK cust_ID
$customers.Cell(4, 2).Shape.TextFrame.TextRange.Text = "Varchar(40)"  # First_name
$customers.Cell(5, 2).Shape.TextFrame.TextRange.Text = "Varchar(40)"  # Last_name
$customers.Cell(6, 2).Shape.TextFrame.TextRange.Text = "Varchar(40)"  # Username

# --- Table 8 ("Orders Table") -----------------------------------------
$orders = $s.Shapes.Item(2).Table
$orders.Cell(3, 2).Shape.TextFrame.TextRange.Text = "Int"       # PK Order_ID
$orders.Cell(4, 2).Shape.TextFrame.TextRange.Text = "Int"       # FK cust ID
$orders.Cell(5, 2).Shape.TextFrame.TextRange.Text = "Dec(7,2)"  # Totalcost
$orders.Rows.Item(6).Delete()                                   # remove FK address_ID row

# --- Table 10 ("Items table") ------------------------------------------
$items = $s.Shapes.Item(3).Table
$items.Cell(3, 2).Shape.TextFrame.TextRange.Text = "Int"       # PK item_ID
$items.Cell(4, 2).Shape.TextFrame.TextRange.Text = "Int"       # Quantity_in_stock
$items.Cell(5, 2).Shape.TextFrame.TextRange.Text = "Dec(7,2)"  # item_price

# --- Table 14 ("Order details") -----------------------------------------
$details = $s.Shapes.Item(4).Table
$details.Cell(3, 2).Shape.TextFrame.TextRange.Text = "Int"       # FK item_ID
$details.Cell(4, 2).Shape.TextFrame.TextRange.Text = "Int"       # FK_orderID
$details.Cell(5, 2).Shape.TextFrame.TextRange.Text = "Dec(7,2)"  # Total_price
$details.Cell(6, 2).Shape.TextFrame.TextRange.Text = "Int"       # quantity
